# B6-PowerPoint.pptx edit: 23 Apr 2020
#
# 1) Re-style the three summary tables (slides 14-16) from the deck's
#    custom "Table_0" style to the built-in table style
#    {41443418-4321-4EA2-BE25-B21C2A941804}.
# 2) Re-colour the deck's theme (theme1.xml, used by the slide master /
#    every slide) from the "Integral / Red Violet" palette to the
#    standard Office palette - font scheme and format scheme are already
#    identical between the two themes shipped in this file, only the
#    12 scheme colours differ.

$p = $ppt.ActivePresentation

# --- 1) Tables: apply the new built-in table style -----------------------
$newTableStyleId = "{41443418-4321-4EA2-BE25-B21C2A941804}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2) Theme colours: Integral/Red Violet -> Office ----------------------
$officeColors = @{
    1  = "000000"  # dk1
    2  = "FFFFFF"  # lt1
    3  = "44546A"  # dk2
    4  = "E7E6E6"  # lt2
    5  = "5B9BD5"  # accent1
    6  = "ED7D31"  # accent2
    7  = "A5A5A5"  # accent3
    8  = "FFC000"  # accent4
    9  = "4472C4"  # accent5
    10 = "70AD47"  # accent6
    11 = "0563C1"  # hlink
    12 = "954F72"  # folHlink
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme
foreach ($idx in $officeColors.Keys) {
    $hex = $officeColors[$idx]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $themeColors.Colors($idx).RGB = $r + ($g * 256) + ($b * 65536)
}
